$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.Value = "'" + $value
    $range.Style = $origStyle
}

Set-TextValue "D2" '36.687.72'
Set-TextValue "E2" '  -0.77%  '
Set-TextValue "D3" '2.060.51'
Set-TextValue "E4" '  -0.10%  '
Set-TextValue "D5" '244.52'
Set-TextValue "E5" '  -0.43%  '
Set-TextValue "D6" '0.667'
Set-TextValue "E6" '  +1.07%  '
Set-TextValue "E7" '  +0.00%  '
Set-TextValue "D8" '55.39'
Set-TextValue "E8" '  -5.02%  '
Set-TextValue "D9" '60.04'
Set-TextValue "E9" '  +1.23%  '
Set-TextValue "D10" '0.368'
Set-TextValue "E10" '  -2.34%  '
Set-TextValue "D11" '0.0751'
Set-TextValue "E11" '  -2.54%  '
Set-TextValue "E12" '  -3.19%  '
Set-TextValue "D13" '0.940'
Set-TextValue "E13" '  +5.50%  '
Set-TextValue "D14" '14.85'
Set-TextValue "E14" '  -3.27%  '
Set-TextValue "D15" '2.360.41'
Set-TextValue "E15" '  +0.49%  '
Set-TextValue "D16" '5.49'
Set-TextValue "E16" '  -3.78%  '
Set-TextValue "D17" '2.059.78'
Set-TextValue "E17" '  +3.08%  '
Set-TextValue "D18" '36.600.18'
Set-TextValue "E18" '  -0.89%  '
Set-TextValue "D19" '17.29'
Set-TextValue "E19" '  -4.44%  '
Set-TextValue "D20" '72.18'
Set-TextValue "E20" '  -2.31%  '
Set-TextValue "E21" '  -1.89%  '
Set-TextValue "D22" '239.16'
Set-TextValue "E22" '  +0.36%  '
Set-TextValue "E23" '  -2.91%  '
Set-TextValue "D24" '0.999'
Set-TextValue "E24" '  -0.11%  '
Set-TextValue "E25" '  -2.27%  '
Set-TextValue "D26" '2.23'
Set-TextValue "E26" '  +4.38%  '
Set-TextValue "D27" '9.28'
Set-TextValue "E27" '  -4.71%  '
Set-TextValue "D28" '165.57'
Set-TextValue "E28" '  -1.76%  '
Set-TextValue "D29" '20.19'
Set-TextValue "E29" '  +0.93%  '
Set-TextValue "E30" '  -0.89%  '
Set-TextValue "E31" '  -7.58%  '
Set-TextValue "E32" '  +6.66%  '
Set-TextValue "D33" '4.50'
Set-TextValue "E33" '  -3.60%  '
Set-TextValue "D34" '0.0599'
Set-TextValue "E34" '  -2.26%  '
Set-TextValue "E35" '  +0.08%  '
Set-TextValue "E36" '  +0.04%  '
Set-TextValue "D37" '0.0844'
Set-TextValue "E37" '  -1.93%  '
Set-TextValue "E38" '  -1.04%  '
Set-TextValue "E39" '  -4.41%  '
Set-TextValue "D40" '5.03'
Set-TextValue "E40" '  -3.79%  '
Set-TextValue "D41" '2.93'
Set-TextValue "E41" '  -4.61%  '
Set-TextValue "E42" '  -2.99%  '
Set-TextValue "E43" '  -3.25%  '
Set-TextValue "D44" '94.93'
Set-TextValue "E44" '  -3.15%  '
Set-TextValue "B45" 'Maker'
Set-TextValue "C45" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D45" '1.410.17'
Set-TextValue "E45" '  +8.54%  '
Set-TextValue "B46" 'Cronos'
Set-TextValue "C46" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D46" '0.0904'
Set-TextValue "E46" '  -6.55%  '
Set-TextValue "B47" 'InjectiveProtocol'
Set-TextValue "C47" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D47" '16.16'
Set-TextValue "E47" '  -4.54%  '
Set-TextValue "B48" 'FraxShare'
Set-TextValue "C48" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D48" '7.57'
Set-TextValue "E48" '  +11.82%  '
Set-TextValue "E49" '  +1.40%  '
Set-TextValue "D50" '2.29'
Set-TextValue "E50" '  -3.34%  '
Set-TextValue "D51" '2.247.90'
Set-TextValue "E51" '  +0.63%  '
